$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The shock ("wstrzas") results for rows 2-9 move from column E into column G,
# carrying their original (right-aligned, text) formatting with them; the
# stale W2..W4 numbers that used to sit in H:J for these rows are cleared out.
for ($r = 2; $r -le 9; $r++) {
    $eCell = $ws.Cells.Item($r, 5)   # column E
    $gCell = $ws.Cells.Item($r, 7)   # column G

    $eCell.Copy($gCell) | Out-Null
    $ws.Range($ws.Cells.Item($r, 8), $ws.Cells.Item($r, 10)).Clear() | Out-Null
    $eCell.ClearContents() | Out-Null
}

# The refreshed UI leaves the new shock-selection column highlighted.
$ws.Range("E2:E9").Select()
